$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"32.50235"
$ws.Range("H2").Value = [double]"97.50704999999999"
$ws.Range("I2").Value = [double]"0.004318312013857221"
$ws.Range("J2").Value = [double]"0.004318312013857221"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.569028"
$ws.Range("N2").Value = [double]"1.707084"
$ws.Range("O2").Value = [double]"0.1016535000995941"
$ws.Range("P2").Value = [double]"0.1016535000995941"
$ws.Range("Q2").Value = [double]"18.4947472158"
$ws.Range("R2").Value = [double]"166.4527249422"
$ws.Range("S2").Value = [double]"0.0004389715307307134"
$ws.Range("T2").Value = [double]"0.0004389715307307134"

$ws.Range("G3").Value = [double]"32.50235"
$ws.Range("H3").Value = [double]"97.50704999999999"
$ws.Range("I3").Value = [double]"0.004318312013857221"
$ws.Range("J3").Value = [double]"0.004318312013857221"
$ws.Range("O3").Value = [double]"0.1962512724671019"
$ws.Range("P3").Value = [double]"0.1962512724671019"
$ws.Range("Q3").Value = [double]"35.705781616"
$ws.Range("R3").Value = [double]"321.352034544"
$ws.Range("S3").Value = [double]"0.0008474742276294532"
$ws.Range("T3").Value = [double]"0.0008474742276294531"

$ws.Range("G4").Value = [double]"32.50235"
$ws.Range("H4").Value = [double]"97.50704999999999"
$ws.Range("I4").Value = [double]"0.004318312013857221"
$ws.Range("J4").Value = [double]"0.004318312013857221"
$ws.Range("M4").Value = [double]"3.083549"
$ws.Range("N4").Value = [double]"9.250647000000001"
$ws.Range("O4").Value = [double]"0.5508578638987945"
$ws.Range("P4").Value = [double]"0.5508578638987945"
$ws.Range("Q4").Value = [double]"100.22258884015"
$ws.Range("R4").Value = [double]"902.00329956135"
$ws.Range("S4").Value = [double]"0.00237877613160189"
$ws.Range("T4").Value = [double]"0.00237877613160189"

$ws.Range("G5").Value = [double]"32.50235"
$ws.Range("H5").Value = [double]"97.50704999999999"
$ws.Range("I5").Value = [double]"0.004318312013857221"
$ws.Range("J5").Value = [double]"0.004318312013857221"
$ws.Range("M5").Value = [double]"0.3400753333333333"
$ws.Range("N5").Value = [double]"1.020226"
$ws.Range("O5").Value = [double]"0.06075245494223394"
$ws.Range("P5").Value = [double]"0.06075245494223393"
$ws.Range("Q5").Value = [double]"11.05324751036667"
$ws.Range("R5").Value = [double]"99.4792275933"
$ws.Range("S5").Value = [double]"0.0002623480560483683"
$ws.Range("T5").Value = [double]"0.0002623480560483683"

$ws.Range("G6").Value = [double]"32.50235"
$ws.Range("H6").Value = [double]"97.50704999999999"
$ws.Range("I6").Value = [double]"0.004318312013857221"
$ws.Range("J6").Value = [double]"0.004318312013857221"
$ws.Range("M6").Value = [double]"0.5065093333333334"
$ws.Range("N6").Value = [double]"1.519528"
$ws.Range("O6").Value = [double]"0.0904849085922755"
$ws.Range("P6").Value = [double]"0.09048490859227548"
$ws.Range("Q6").Value = [double]"16.46274363026667"
$ws.Range("R6").Value = [double]"148.1646926724"
$ws.Range("S6").Value = [double]"0.0003907420678467958"
$ws.Range("T6").Value = [double]"0.0003907420678467957"

$ws.Range("I7").Value = [double]"0.006762540683959845"
$ws.Range("J7").Value = [double]"0.006762540683959845"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"0.569028"
$ws.Range("N7").Value = [double]"1.707084"
$ws.Range("O7").Value = [double]"0.1016535000995941"
$ws.Range("P7").Value = [double]"0.1016535000995941"
$ws.Range("Q7").Value = [double]"28.963048544212"
$ws.Range("R7").Value = [double]"260.667436897908"
$ws.Range("S7").Value = [double]"0.0006874359300904213"
$ws.Range("T7").Value = [double]"0.0006874359300904212"

$ws.Range("I8").Value = [double]"0.006762540683959845"
$ws.Range("J8").Value = [double]"0.006762540683959845"
$ws.Range("O8").Value = [double]"0.1962512724671019"
$ws.Range("P8").Value = [double]"0.1962512724671019"
$ws.Range("Q8").Value = [double]"55.91578377290666"
$ws.Range("S8").Value = [double]"0.001327157214337665"
$ws.Range("T8").Value = [double]"0.001327157214337665"

$ws.Range("I9").Value = [double]"0.006762540683959845"
$ws.Range("J9").Value = [double]"0.006762540683959845"
$ws.Range("M9").Value = [double]"3.083549"
$ws.Range("N9").Value = [double]"9.250647000000001"
$ws.Range("O9").Value = [double]"0.5508578638987945"
$ws.Range("P9").Value = [double]"0.5508578638987945"
$ws.Range("Q9").Value = [double]"156.9500611137877"
$ws.Range("R9").Value = [double]"1412.550550024089"
$ws.Range("S9").Value = [double]"0.003725198715694814"
$ws.Range("T9").Value = [double]"0.003725198715694814"

$ws.Range("I10").Value = [double]"0.006762540683959845"
$ws.Range("J10").Value = [double]"0.006762540683959845"
$ws.Range("M10").Value = [double]"0.3400753333333333"
$ws.Range("N10").Value = [double]"1.020226"
$ws.Range("O10").Value = [double]"0.06075245494223394"
$ws.Range("P10").Value = [double]"0.06075245494223393"
$ws.Range("Q10").Value = [double]"17.30954959689578"
$ws.Range("R10").Value = [double]"155.785946372062"
$ws.Range("S10").Value = [double]"0.0004108409481972944"
$ws.Range("T10").Value = [double]"0.0004108409481972944"

$ws.Range("I11").Value = [double]"0.006762540683959845"
$ws.Range("J11").Value = [double]"0.006762540683959845"
$ws.Range("M11").Value = [double]"0.5065093333333334"
$ws.Range("N11").Value = [double]"1.519528"
$ws.Range("O11").Value = [double]"0.0904849085922755"
$ws.Range("P11").Value = [double]"0.09048490859227548"
$ws.Range("Q11").Value = [double]"25.78090078068178"
$ws.Range("R11").Value = [double]"232.028107026136"
$ws.Range("S11").Value = [double]"0.0006119078756396508"
$ws.Range("T11").Value = [double]"0.0006119078756396507"

$ws.Range("G12").Value = [double]"3274.382486666667"
$ws.Range("H12").Value = [double]"9823.14746"
$ws.Range("I12").Value = [double]"0.4350394734576531"
$ws.Range("J12").Value = [double]"0.435039473457653"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"0.569028"
$ws.Range("N12").Value = [double]"1.707084"
$ws.Range("O12").Value = [double]"0.1016535000995941"
$ws.Range("P12").Value = [double]"0.1016535000995941"
$ws.Range("Q12").Value = [double]"1863.21531762296"
$ws.Range("R12").Value = [double]"16768.93785860664"
$ws.Range("S12").Value = [double]"0.0442232851584549"
$ws.Range("T12").Value = [double]"0.04422328515845489"

$ws.Range("G13").Value = [double]"3274.382486666667"
$ws.Range("H13").Value = [double]"9823.14746"
$ws.Range("I13").Value = [double]"0.4350394734576531"
$ws.Range("J13").Value = [double]"0.435039473457653"
$ws.Range("O13").Value = [double]"0.1962512724671019"
$ws.Range("P13").Value = [double]"0.1962512724671019"
$ws.Range("Q13").Value = [double]"3597.105624552533"
$ws.Range("R13").Value = [double]"32373.9506209728"
$ws.Range("S13").Value = [double]"0.08537705023948244"
$ws.Range("T13").Value = [double]"0.08537705023948242"

$ws.Range("G14").Value = [double]"3274.382486666667"
$ws.Range("H14").Value = [double]"9823.14746"
$ws.Range("I14").Value = [double]"0.4350394734576531"
$ws.Range("J14").Value = [double]"0.435039473457653"
$ws.Range("M14").Value = [double]"3.083549"
$ws.Range("N14").Value = [double]"9.250647000000001"
$ws.Range("O14").Value = [double]"0.5508578638987945"
$ws.Range("P14").Value = [double]"0.5508578638987945"
$ws.Range("Q14").Value = [double]"10096.71884237851"
$ws.Range("R14").Value = [double]"90870.46958140662"
$ws.Range("S14").Value = [double]"0.2396449150605391"
$ws.Range("T14").Value = [double]"0.2396449150605391"

$ws.Range("G15").Value = [double]"3274.382486666667"
$ws.Range("H15").Value = [double]"9823.14746"
$ws.Range("I15").Value = [double]"0.4350394734576531"
$ws.Range("J15").Value = [double]"0.435039473457653"
$ws.Range("M15").Value = [double]"0.3400753333333333"
$ws.Range("N15").Value = [double]"1.020226"
$ws.Range("O15").Value = [double]"0.06075245494223394"
$ws.Range("P15").Value = [double]"0.06075245494223393"
$ws.Range("Q15").Value = [double]"1113.536715613996"
$ws.Range("R15").Value = [double]"10021.83044052596"
$ws.Range("S15").Value = [double]"0.02642971600932925"
$ws.Range("T15").Value = [double]"0.02642971600932924"

$ws.Range("G16").Value = [double]"3274.382486666667"
$ws.Range("H16").Value = [double]"9823.14746"
$ws.Range("I16").Value = [double]"0.4350394734576531"
$ws.Range("J16").Value = [double]"0.435039473457653"
$ws.Range("M16").Value = [double]"0.5065093333333334"
$ws.Range("N16").Value = [double]"1.519528"
$ws.Range("O16").Value = [double]"0.0904849085922755"
$ws.Range("P16").Value = [double]"0.09048490859227548"
$ws.Range("Q16").Value = [double]"1658.505290399876"
$ws.Range("R16").Value = [double]"14926.54761359888"
$ws.Range("S16").Value = [double]"0.0393645069898474"
$ws.Range("T16").Value = [double]"0.03936450698984739"

$ws.Range("G17").Value = [double]"7.278837333333333"
$ws.Range("H17").Value = [double]"21.836512"
$ws.Range("I17").Value = [double]"0.0009670774791190726"
$ws.Range("J17").Value = [double]"0.0009670774791190726"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.569028"
$ws.Range("N17").Value = [double]"1.707084"
$ws.Range("O17").Value = [double]"0.1016535000995941"
$ws.Range("P17").Value = [double]"0.1016535000995941"
$ws.Range("Q17").Value = [double]"4.141862250111999"
$ws.Range("R17").Value = [double]"37.276760251008"
$ws.Range("S17").Value = [double]"9.830681061994586E-05"
$ws.Range("T17").Value = [double]"9.830681061994585E-05"

$ws.Range("G18").Value = [double]"7.278837333333333"
$ws.Range("H18").Value = [double]"21.836512"
$ws.Range("I18").Value = [double]"0.0009670774791190726"
$ws.Range("J18").Value = [double]"0.0009670774791190726"
$ws.Range("O18").Value = [double]"0.1962512724671019"
$ws.Range("P18").Value = [double]"0.1962512724671019"
$ws.Range("Q18").Value = [double]"7.996239540906666"
$ws.Range("R18").Value = [double]"71.96615586816"
$ws.Range("S18").Value = [double]"0.0001897901858513952"
$ws.Range("T18").Value = [double]"0.0001897901858513952"

$ws.Range("G19").Value = [double]"7.278837333333333"
$ws.Range("H19").Value = [double]"21.836512"
$ws.Range("I19").Value = [double]"0.0009670774791190726"
$ws.Range("J19").Value = [double]"0.0009670774791190726"
$ws.Range("M19").Value = [double]"3.083549"
$ws.Range("N19").Value = [double]"9.250647000000001"
$ws.Range("O19").Value = [double]"0.5508578638987945"
$ws.Range("P19").Value = [double]"0.5508578638987945"
$ws.Range("Q19").Value = [double]"22.44465158036267"
$ws.Range("R19").Value = [double]"202.001864223264"
$ws.Range("S19").Value = [double]"0.0005327222343721634"
$ws.Range("T19").Value = [double]"0.0005327222343721634"

$ws.Range("G20").Value = [double]"7.278837333333333"
$ws.Range("H20").Value = [double]"21.836512"
$ws.Range("I20").Value = [double]"0.0009670774791190726"
$ws.Range("J20").Value = [double]"0.0009670774791190726"
$ws.Range("M20").Value = [double]"0.3400753333333333"
$ws.Range("N20").Value = [double]"1.020226"
$ws.Range("O20").Value = [double]"0.06075245494223394"
$ws.Range("P20").Value = [double]"0.06075245494223393"
$ws.Range("Q20").Value = [double]"2.475353032412444"
$ws.Range("R20").Value = [double]"22.278177291712"
$ws.Range("S20").Value = [double]"5.875233097583065E-05"
$ws.Range("T20").Value = [double]"5.875233097583064E-05"

$ws.Range("G21").Value = [double]"7.278837333333333"
$ws.Range("H21").Value = [double]"21.836512"
$ws.Range("I21").Value = [double]"0.0009670774791190726"
$ws.Range("J21").Value = [double]"0.0009670774791190726"
$ws.Range("M21").Value = [double]"0.5065093333333334"
$ws.Range("N21").Value = [double]"1.519528"
$ws.Range("O21").Value = [double]"0.0904849085922755"
$ws.Range("P21").Value = [double]"0.09048490859227548"
$ws.Range("Q21").Value = [double]"3.686799045148444"
$ws.Range("R21").Value = [double]"33.181191406336"
$ws.Range("S21").Value = [double]"8.750591729973751E-05"
$ws.Range("T21").Value = [double]"8.750591729973749E-05"

$ws.Range("G22").Value = [double]"4161.570231333333"
$ws.Range("H22").Value = [double]"12484.710694"
$ws.Range("I22").Value = [double]"0.5529125963654108"
$ws.Range("J22").Value = [double]"0.5529125963654108"
$ws.Range("K22").Value = [double]"3"
$ws.Range("L22").Value = [double]"1"
$ws.Range("M22").Value = [double]"0.569028"
$ws.Range("N22").Value = [double]"1.707084"
$ws.Range("O22").Value = [double]"0.1016535000995941"
$ws.Range("P22").Value = [double]"0.1016535000995941"
$ws.Range("Q22").Value = [double]"2368.049985595144"
$ws.Range("R22").Value = [double]"21312.4498703563"
$ws.Range("S22").Value = [double]"0.05620550066969811"
$ws.Range("T22").Value = [double]"0.05620550066969811"

$ws.Range("G23").Value = [double]"4161.570231333333"
$ws.Range("H23").Value = [double]"12484.710694"
$ws.Range("I23").Value = [double]"0.5529125963654108"
$ws.Range("J23").Value = [double]"0.5529125963654108"
$ws.Range("O23").Value = [double]"0.1962512724671019"
$ws.Range("P23").Value = [double]"0.1962512724671019"
$ws.Range("Q23").Value = [double]"4571.734593333546"
$ws.Range("R23").Value = [double]"41145.61134000192"
$ws.Range("S23").Value = [double]"0.108509800599801"
$ws.Range("T23").Value = [double]"0.108509800599801"

$ws.Range("G24").Value = [double]"4161.570231333333"
$ws.Range("H24").Value = [double]"12484.710694"
$ws.Range("I24").Value = [double]"0.5529125963654108"
$ws.Range("J24").Value = [double]"0.5529125963654108"
$ws.Range("M24").Value = [double]"3.083549"
$ws.Range("N24").Value = [double]"9.250647000000001"
$ws.Range("O24").Value = [double]"0.5508578638987945"
$ws.Range("P24").Value = [double]"0.5508578638987945"
$ws.Range("Q24").Value = [double]"12832.40572525767"
$ws.Range("R24").Value = [double]"115491.651527319"
$ws.Range("S24").Value = [double]"0.3045762517565866"
$ws.Range("T24").Value = [double]"0.3045762517565866"

$ws.Range("G25").Value = [double]"4161.570231333333"
$ws.Range("H25").Value = [double]"12484.710694"
$ws.Range("I25").Value = [double]"0.5529125963654108"
$ws.Range("J25").Value = [double]"0.5529125963654108"
$ws.Range("M25").Value = [double]"0.3400753333333333"
$ws.Range("N25").Value = [double]"1.020226"
$ws.Range("O25").Value = [double]"0.06075245494223394"
$ws.Range("P25").Value = [double]"0.06075245494223393"
$ws.Range("Q25").Value = [double]"1415.24738361076"
$ws.Range("R25").Value = [double]"12737.22645249685"
$ws.Range("S25").Value = [double]"0.0335907975976832"
$ws.Range("T25").Value = [double]"0.0335907975976832"

$ws.Range("G26").Value = [double]"4161.570231333333"
$ws.Range("H26").Value = [double]"12484.710694"
$ws.Range("I26").Value = [double]"0.5529125963654108"
$ws.Range("J26").Value = [double]"0.5529125963654108"
$ws.Range("M26").Value = [double]"0.5065093333333334"
$ws.Range("N26").Value = [double]"1.519528"
$ws.Range("O26").Value = [double]"0.0904849085922755"
$ws.Range("P26").Value = [double]"0.09048490859227548"
$ws.Range("Q26").Value = [double]"2107.874163492493"
$ws.Range("R26").Value = [double]"18970.86747143243"
$ws.Range("S26").Value = [double]"0.05003024574164191"
$ws.Range("T26").Value = [double]"0.0500302457416419"
